$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated values (new TPM) for existing data rows 2-21
$ws.Range("G2").Value = 0.052224
$ws.Range("H2").Value = 0.156672
$ws.Range("I2").Value = 0.04490044548800991
$ws.Range("J2").Value = 0.05567964586107596
$ws.Range("M2").Value = 162.7225033333333
$ws.Range("N2").Value = 488.16751
$ws.Range("O2").Value = 0.5231437953541009
$ws.Range("P2").Value = 0.5247717033381212
$ws.Range("Q2").Value = 8.49802001408
$ws.Range("R2").Value = 76.48218012672
$ws.Range("S2").Value = 0.02348938946568742
$ws.Range("T2").Value = 0.0292191025997802

$ws.Range("G3").Value = 0.052224
$ws.Range("H3").Value = 0.156672
$ws.Range("I3").Value = 0.04490044548800991
$ws.Range("J3").Value = 0.05567964586107596
$ws.Range("O3").Value = 0.0009322191998643353
$ws.Range("P3").Value = 0.0009351200601857102
$ws.Range("Q3").Value = 0.015143097344
$ws.Range("R3").Value = 0.136287876096
$ws.Range("S3").Value = 0.0000418570573663848
$ws.Range("T3").Value = 0.00005206715378872838

$ws.Range("G4").Value = 0.052224
$ws.Range("H4").Value = 0.156672
$ws.Range("I4").Value = 0.04490044548800991
$ws.Range("J4").Value = 0.05567964586107596
$ws.Range("M4").Value = 61.580654
$ws.Range("N4").Value = 184.741962
$ws.Range("O4").Value = 0.1979783766474813
$ws.Range("P4").Value = 0.1985944416431287
$ws.Range("Q4").Value = 3.215988074496
$ws.Range("R4").Value = 28.943892670464
$ws.Range("S4").Value = 0.008889317308464928
$ws.Range("T4").Value = 0.01105766818066752

$ws.Range("G5").Value = 0.052224
$ws.Range("H5").Value = 0.156672
$ws.Range("I5").Value = 0.04490044548800991
$ws.Range("J5").Value = 0.05567964586107596
$ws.Range("M5").Value = 2.8947245
$ws.Range("N5").Value = 5.789449
$ws.Range("O5").Value = 0.009306378223129816
$ws.Range("P5").Value = 0.00622355841157717
$ws.Range("Q5").Value = 0.151174092288
$ws.Range("R5").Value = 0.9070445537280001
$ws.Range("S5").Value = 0.0004178605280984428
$ws.Range("T5").Value = 0.0003465255283523373

$ws.Range("G6").Value = 0.052224
$ws.Range("H6").Value = 0.156672
$ws.Range("I6").Value = 0.04490044548800991
$ws.Range("J6").Value = 0.05567964586107596
$ws.Range("M6").Value = 83.559527
$ws.Range("N6").Value = 250.678581
$ws.Range("O6").Value = 0.2686392305754237
$ws.Range("P6").Value = 0.2694751765469873
$ws.Range("Q6").Value = 4.363812738048
$ws.Range("R6").Value = 39.274314642432
$ws.Range("S6").Value = 0.01206202112839273
$ws.Range("T6").Value = 0.01500428239848717

$ws.Range("I7").Value = 0.3275750472262571
$ws.Range("J7").Value = 0.40621562713345
$ws.Range("M7").Value = 162.7225033333333
$ws.Range("N7").Value = 488.16751
$ws.Range("O7").Value = 0.5231437953541009
$ws.Range("P7").Value = 0.5247717033381212
$ws.Range("Q7").Value = 61.99803314168221
$ws.Range("R7").Value = 557.98229827514
$ws.Range("S7").Value = 0.171368853469243
$ws.Range("T7").Value = 0.2131704665733837

$ws.Range("I8").Value = 0.3275750472262571
$ws.Range("J8").Value = 0.40621562713345
$ws.Range("O8").Value = 0.0009322191998643353
$ws.Range("P8").Value = 0.0009351200601857102
$ws.Range("S8").Value = 0.0003053717484207833
$ws.Range("T8").Value = 0.0003798603816934078

$ws.Range("I9").Value = 0.3275750472262571
$ws.Range("J9").Value = 0.40621562713345
$ws.Range("M9").Value = 61.580654
$ws.Range("N9").Value = 184.741962
$ws.Range("O9").Value = 0.1979783766474813
$ws.Range("P9").Value = 0.1985944416431287
$ws.Range("Q9").Value = 23.46251655038533
$ws.Range("R9").Value = 211.162648953468
$ws.Range("S9").Value = 0.0648527760800764
$ws.Range("T9").Value = 0.08067216565728086

$ws.Range("I10").Value = 0.3275750472262571
$ws.Range("J10").Value = 0.40621562713345
$ws.Range("M10").Value = 2.8947245
$ws.Range("N10").Value = 5.789449
$ws.Range("O10").Value = 0.009306378223129816
$ws.Range("P10").Value = 0.00622355841157717
$ws.Range("Q10").Value = 1.102903543214333
$ws.Range("R10").Value = 6.617421259286
$ws.Range("S10").Value = 0.00304853728594716
$ws.Range("T10").Value = 0.002528106683160478

$ws.Range("I11").Value = 0.3275750472262571
$ws.Range("J11").Value = 0.40621562713345
$ws.Range("M11").Value = 83.559527
$ws.Range("N11").Value = 250.678581
$ws.Range("O11").Value = 0.2686392305754237
$ws.Range("P11").Value = 0.2694751765469873
$ws.Range("Q11").Value = 31.83656973145933
$ws.Range("R11").Value = 286.529127583134
$ws.Range("S11").Value = 0.08799950864256978
$ws.Range("T11").Value = 0.1094650278379316

$ws.Range("G12").Value = 0.01751533333333333
$ws.Range("H12").Value = 0.052546
$ws.Range("I12").Value = 0.01505909676657583
$ws.Range("J12").Value = 0.01867431750035806
$ws.Range("M12").Value = 162.7225033333333
$ws.Range("N12").Value = 488.16751
$ws.Range("O12").Value = 0.5231437953541009
$ws.Range("P12").Value = 0.5247717033381212
$ws.Range("Q12").Value = 2.850138886717778
$ws.Range("R12").Value = 25.65124998046
$ws.Range("S12").Value = 0.007878073037071151
$ws.Range("T12").Value = 0.009799753403339783

$ws.Range("G13").Value = 0.01751533333333333
$ws.Range("H13").Value = 0.052546
$ws.Range("I13").Value = 0.01505909676657583
$ws.Range("J13").Value = 0.01867431750035806
$ws.Range("O13").Value = 0.0009322191998643353
$ws.Range("P13").Value = 0.0009351200601857102
$ws.Range("Q13").Value = 0.005078821953111112
$ws.Range("R13").Value = 0.045709397578
$ws.Range("S13").Value = 0.00001403837913841692
$ws.Range("T13").Value = 0.00001746272890486189

$ws.Range("G14").Value = 0.01751533333333333
$ws.Range("H14").Value = 0.052546
$ws.Range("I14").Value = 0.01505909676657583
$ws.Range("J14").Value = 0.01867431750035806
$ws.Range("M14").Value = 61.580654
$ws.Range("N14").Value = 184.741962
$ws.Range("O14").Value = 0.1979783766474813
$ws.Range("P14").Value = 0.1985944416431287
$ws.Range("Q14").Value = 1.078605681694667
$ws.Range("R14").Value = 9.707451135252
$ws.Range("S14").Value = 0.002981375531624018
$ws.Range("T14").Value = 0.003708615657050115

$ws.Range("G15").Value = 0.01751533333333333
$ws.Range("H15").Value = 0.052546
$ws.Range("I15").Value = 0.01505909676657583
$ws.Range("J15").Value = 0.01867431750035806
$ws.Range("M15").Value = 2.8947245
$ws.Range("N15").Value = 5.789449
$ws.Range("O15").Value = 0.009306378223129816
$ws.Range("P15").Value = 0.00622355841157717
$ws.Range("Q15").Value = 0.05070206452566667
$ws.Range("R15").Value = 0.304212387154
$ws.Range("S15").Value = 0.0001401456502084659
$ws.Range("T15").Value = 0.0001162207057598161

$ws.Range("G16").Value = 0.01751533333333333
$ws.Range("H16").Value = 0.052546
$ws.Range("I16").Value = 0.01505909676657583
$ws.Range("J16").Value = 0.01867431750035806
$ws.Range("M16").Value = 83.559527
$ws.Range("N16").Value = 250.678581
$ws.Range("O16").Value = 0.2686392305754237
$ws.Range("P16").Value = 0.2694751765469873
$ws.Range("Q16").Value = 1.463572968580667
$ws.Range("R16").Value = 13.172156717226
$ws.Range("S16").Value = 0.004045464168533782
$ws.Range("T16").Value = 0.005032265005303481

$ws.Range("G17").Value = 0.6755085000000001
$ws.Range("H17").Value = 1.351017
$ws.Range("I17").Value = 0.5807795760749338
$ws.Range("J17").Value = 0.4801377917706627
$ws.Range("M17").Value = 162.7225033333333
$ws.Range("N17").Value = 488.16751
$ws.Range("O17").Value = 0.5231437953541009
$ws.Range("P17").Value = 0.5247717033381212
$ws.Range("Q17").Value = 109.920434142945
$ws.Range("R17").Value = 659.5226048576701
$ws.Range("S17").Value = 0.3038312316919867
$ws.Range("T17").Value = 0.2519627268244948

$ws.Range("G18").Value = 0.6755085000000001
$ws.Range("H18").Value = 1.351017
$ws.Range("I18").Value = 0.5807795760749338
$ws.Range("J18").Value = 0.4801377917706627
$ws.Range("O18").Value = 0.0009322191998643353
$ws.Range("P18").Value = 0.0009351200601857102
$ws.Range("Q18").Value = 0.1958733718635
$ws.Range("R18").Value = 1.175240231181
$ws.Range("S18").Value = 0.0005414138717061227
$ws.Range("T18").Value = 0.0004489864807380161

$ws.Range("G19").Value = 0.6755085000000001
$ws.Range("H19").Value = 1.351017
$ws.Range("I19").Value = 0.5807795760749338
$ws.Range("J19").Value = 0.4801377917706627
$ws.Range("M19").Value = 61.580654
$ws.Range("N19").Value = 184.741962
$ws.Range("O19").Value = 0.1979783766474813
$ws.Range("P19").Value = 0.1985944416431287
$ws.Range("Q19").Value = 41.598255212559
$ws.Range("R19").Value = 249.589531275354
$ws.Range("S19").Value = 0.1149817976613278
$ws.Range("T19").Value = 0.09535269666845957

$ws.Range("G20").Value = 0.6755085000000001
$ws.Range("H20").Value = 1.351017
$ws.Range("I20").Value = 0.5807795760749338
$ws.Range("J20").Value = 0.4801377917706627
$ws.Range("M20").Value = 2.8947245
$ws.Range("N20").Value = 5.789449
$ws.Range("O20").Value = 0.009306378223129816
$ws.Range("P20").Value = 0.00622355841157717
$ws.Range("Q20").Value = 1.95541100490825
$ws.Range("R20").Value = 7.821644019633001
$ws.Range("S20").Value = 0.00540495439922233
$ws.Range("T20").Value = 0.002988165592690396

$ws.Range("G21").Value = 0.6755085000000001
$ws.Range("H21").Value = 1.351017
$ws.Range("I21").Value = 0.5807795760749338
$ws.Range("J21").Value = 0.4801377917706627
$ws.Range("M21").Value = 83.559527
$ws.Range("N21").Value = 250.678581
$ws.Range("O21").Value = 0.2686392305754237
$ws.Range("P21").Value = 0.2694751765469873
$ws.Range("Q21").Value = 56.4451707444795
$ws.Range("R21").Value = 338.671024466877
$ws.Range("S21").Value = 0.1560201784506909
$ws.Range("T21").Value = 0.12938521620428

# Append new rows 22-26 (Resolving-Mac as sending cluster)
$ws.Range("A22").Value = "Resolving-Mac"
$ws.Range("B22").Value = "Col4a6"
$ws.Range("C22").Value = "Cd93"
$ws.Range("D22").Value = "ECs"
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 0.3333333333333333
$ws.Range("G22").Value = 0.036854
$ws.Range("H22").Value = 0.110562
$ws.Range("I22").Value = 0.03168583444422329
$ws.Range("J22").Value = 0.03929261773445338
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 162.7225033333333
$ws.Range("N22").Value = 488.16751
$ws.Range("O22").Value = 0.5231437953541009
$ws.Range("P22").Value = 0.5247717033381212
$ws.Range("Q22").Value = 5.996975137846666
$ws.Range("R22").Value = 53.97277624061999
$ws.Range("S22").Value = 0.01657624769011267
$ws.Range("T22").Value = 0.02061965393712277

$ws.Range("A23").Value = "Resolving-Mac"
$ws.Range("B23").Value = "Col4a6"
$ws.Range("C23").Value = "Cd93"
$ws.Range("D23").Value = "FAPs"
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 0.3333333333333333
$ws.Range("G23").Value = 0.036854
$ws.Range("H23").Value = 0.110562
$ws.Range("I23").Value = 0.03168583444422329
$ws.Range("J23").Value = 0.03929261773445338
$ws.Range("K23").Value = 1
$ws.Range("L23").Value = 0.3333333333333333
$ws.Range("M23").Value = 0.2899643333333333
$ws.Range("N23").Value = 0.869893
$ws.Range("O23").Value = 0.0009322191998643353
$ws.Range("P23").Value = 0.0009351200601857102
$ws.Range("Q23").Value = 0.01068634554066667
$ws.Range("R23").Value = 0.09617710986599999
$ws.Range("S23").Value = 0.00002953814323262764
$ws.Range("T23").Value = 0.00003674331506069615

$ws.Range("A24").Value = "Resolving-Mac"
$ws.Range("B24").Value = "Col4a6"
$ws.Range("C24").Value = "Cd93"
$ws.Range("D24").Value = "Inflammatory-Mac"
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 0.3333333333333333
$ws.Range("G24").Value = 0.036854
$ws.Range("H24").Value = 0.110562
$ws.Range("I24").Value = 0.03168583444422329
$ws.Range("J24").Value = 0.03929261773445338
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 61.580654
$ws.Range("N24").Value = 184.741962
$ws.Range("O24").Value = 0.1979783766474813
$ws.Range("P24").Value = 0.1985944416431287
$ws.Range("Q24").Value = 2.269493422516
$ws.Range("R24").Value = 20.425440802644
$ws.Range("S24").Value = 0.006273110065988175
$ws.Range("T24").Value = 0.007803295479670666

$ws.Range("A25").Value = "Resolving-Mac"
$ws.Range("B25").Value = "Col4a6"
$ws.Range("C25").Value = "Cd93"
$ws.Range("D25").Value = "MuSCs"
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 0.3333333333333333
$ws.Range("G25").Value = 0.036854
$ws.Range("H25").Value = 0.110562
$ws.Range("I25").Value = 0.03168583444422329
$ws.Range("J25").Value = 0.03929261773445338
$ws.Range("K25").Value = 2
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 2.8947245
$ws.Range("N25").Value = 5.789449
$ws.Range("O25").Value = 0.009306378223129816
$ws.Range("P25").Value = 0.00622355841157717
$ws.Range("Q25").Value = 0.106682176723
$ws.Range("R25").Value = 0.640093060338
$ws.Range("S25").Value = 0.0002948803596534163
$ws.Range("T25").Value = 0.0002445399016141436

$ws.Range("A26").Value = "Resolving-Mac"
$ws.Range("B26").Value = "Col4a6"
$ws.Range("C26").Value = "Cd93"
$ws.Range("D26").Value = "Resolving-Mac"
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 0.3333333333333333
$ws.Range("G26").Value = 0.036854
$ws.Range("H26").Value = 0.110562
$ws.Range("I26").Value = 0.03168583444422329
$ws.Range("J26").Value = 0.03929261773445338
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 83.559527
$ws.Range("N26").Value = 250.678581
$ws.Range("O26").Value = 0.2686392305754237
$ws.Range("P26").Value = 0.2694751765469873
$ws.Range("Q26").Value = 3.079502808058
$ws.Range("R26").Value = 27.715525272522
$ws.Range("S26").Value = 0.008512058185236402
$ws.Range("T26").Value = 0.01058838510098511
